# Update the "Estado de Cuenta" table: previous account-statement periods
# are removed and new ones are added; the underlying data (valor mora /
# salario basico) is refreshed accordingly.
#
# Net effect observed in the target workbook: the period column (E16:E24)
# now lists the periods in reverse chronological order (1807 down to 1607),
# and the "Valor Mora" (F) / "Salario Basico" (G) columns are updated to
# match the new period for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @("1807", "1802", "1705", "1702", "1701", "1611", "1610", "1608", "1607")
$valorMora = @(31249, 29509, 27578, 27578, 27578, 27578, 27578, 27578, 27578)
$salarioBasico = 781242

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    $ws.Cells.Item($row, 6).Value = $valorMora[$i]
    $ws.Cells.Item($row, 7).Value = $salarioBasico
}
